# "add tank color, add hit system"
#
# Slide 1 holds the weekly-progress table ("표 4"). Two cells in the
# "percent complete" (4th) column need updating:
#   - Week 5 (2인 플레이 / hit system)  : 10%  -> 90%
#   - Week 6 (적군 AI)                  : ""   -> 10%
#
# NOTE: Week 5 is updated before Week 6 on purpose - Week 5 already
# contains the literal text "10%", and clearing that value out of the
# table first keeps the newly-typed "10%" in Week 6 from being treated
# as a re-used/duplicate run when PowerPoint re-generates its paragraph
# end-of-run properties.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

$table = $s.Shapes.Item("표 4").Table

# Week 5 row ("2인 플레이") - progress 10% -> 90%
$week5Progress = $table.Cell(5, 4)
$week5Progress.Shape.TextFrame.TextRange.Text = "90%"

# Week 6 row ("적군 AI") - progress was blank, now 10%
$week6Progress = $table.Cell(6, 4)
$week6Progress.Shape.TextFrame.TextRange.Text = "10%"
